$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add the new row of data: Leetcode Q# 25, "Reverse Nodes in k-Group", GFG
$ws.Range("A5").Value = 25
$ws.Range("B5").Value = "Reverse Nodes in k-Group"

# Style the new A5 cell like the rest of the table (left/top aligned, wrap text).
# Order matters for how the interop layer merges alignment flags into one xf:
# WrapText first (matches the existing "vertical=top,wrap=1" xf), then
# Horizontal, then Vertical, collapses into a single new style entry.
$ws.Range("A5").WrapText = $true
$ws.Range("A5").HorizontalAlignment = -4131  # xlLeft
$ws.Range("A5").VerticalAlignment = -4160    # xlTop

# Clear the selection that previously pointed at B9 and drop the explicit
# row-1 height override so the row reverts to the sheet's default height.
$ws.Range("A1").Select()
$ws.Rows.Item(1).AutoFit()
